{"js": "// Append a new bulleted \"Things I Learned\" item to the end of the document,\n// matching the formatting (ListParagraph style, same numbered/bulleted list)\n// of the preceding list item.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body is the final bullet item\n// (\"Instructor used Mathf.Atan2 ...\") that the new item should follow.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Remember which bullet/numbered list the final paragraph belongs to so the\n// new paragraph can be attached to the very same list (same numId).\nconst existingList = lastParagraph.list;\nexistingList.load(\"id\");\nawait context.sync();\n\nconst newText =\n  \"Input.GetAxis Gets a range of values between -1 and 1 and it is based on sensitivity which is why when we press W or A or S or D, the tank moves in that direction but keeps on moving and the solution for this is using Input.GetAxisRaw instead as that only gets values -1 and 1 which will make the tank stop as soon as we let go of the buttons.\";\n\nconst newParagraph = lastParagraph.insertParagraph(newText, Word.InsertLocation.after);\nnewParagraph.style = \"List Paragraph\";\nnewParagraph.attachToList(existingList.id, 0);\n\nawait context.sync();\n", "ps1": "# Append a new bulleted \"Things I Learned\" item to the end of the document,\n# matching the formatting (ListParagraph style, same numbered/bulleted list)\n# of the preceding list item.\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the document is the final bullet item\n# (\"Instructor used Mathf.Atan2 ...\") that the new item should follow.\n$lastParagraph = $d.Paragraphs.Last\n\n# Inserting a paragraph break right after this paragraph's range creates a\n# new paragraph that inherits the same style / list (numId) formatting.\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Input.GetAxis Gets a range of values between -1 and 1 and it is based on sensitivity which is why when we press W or A or S or D, the tank moves in that direction but keeps on moving and the solution for this is using Input.GetAxisRaw instead as that only gets values -1 and 1 which will make the tank stop as soon as we let go of the buttons.\"\n"}
